# Daily auto-push update: insert a new "2026/02/15 19:00" reading just
# before the "2026/12/29" block (row 826) and shift every subsequent
# row down by one (the former last row, row 867, becomes row 868).
#
# Net effect on the sheet:
#   - one row is inserted at row 826
#   - dimension grows from A1:D867 to A1:D868
#   - rows 826..867 (old) become rows 827..868 (new), values unchanged
#   - new row 826 gets: A=2026/02/15  B (weekday)  C=19  D=201

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push everything from row 826 down by one row.
$ws.Rows.Item(826).Insert()

# Columns A/B (date / weekday) for the new row are identical text to the
# row directly above (row 825 is also "2026/02/15"), so copy those two
# cells instead of typing the values in. Copying preserves the plain
# text cell type (matching the rest of the sheet) instead of letting the
# COM layer's "smart" input parser turn "2026/02/15" into a real date
# serial number when assigned via .Value.
$ws.Range("A825:B825").Copy($ws.Range("A826:B826"))

# Numeric columns are safe to set directly.
$ws.Range("C826").Value = 19
$ws.Range("D826").Value = 201
